$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 1736973.5
$ws.Cells.Item(18, 9).Value = 1984969.8
$ws.Cells.Item(18, 11).Value = 1984969.8
$ws.Cells.Item(18, 13).Value = -1984685.8
$ws.Cells.Item(19, 8).Value = 6803200.5
$ws.Cells.Item(19, 9).Value = 12987441
$ws.Cells.Item(19, 10).Value = 535.7
$ws.Cells.Item(19, 11).Value = 12987441
$ws.Cells.Item(19, 12).Value = 535.7
$ws.Cells.Item(19, 13).Value = -12987266
$ws.Cells.Item(19, 14).Value = -885.7
$ws.Cells.Item(76, 8).Value = 2945.6
$ws.Cells.Item(76, 9).Value = 2806.2942
$ws.Cells.Item(76, 11).Value = 2806.2942
$ws.Cells.Item(76, 13).Value = -2491.2942
$ws.Cells.Item(79, 8).Value = 2945.6
$ws.Cells.Item(79, 9).Value = 2806.2942
$ws.Cells.Item(79, 11).Value = 2806.2942
$ws.Cells.Item(79, 13).Value = -1714.2942
$ws.Cells.Item(112, 8).Value = 1340.68
$ws.Cells.Item(112, 10).Value = 1353.7551
$ws.Cells.Item(112, 12).Value = 4061.2653
$ws.Cells.Item(112, 14).Value = -6277.2653
$ws.Cells.Item(113, 8).Value = 3292.4
$ws.Cells.Item(113, 9).Value = 3606.875
$ws.Cells.Item(113, 10).Value = 2733.3333
$ws.Cells.Item(113, 11).Value = 3606.875
$ws.Cells.Item(113, 12).Value = 2733.3333
$ws.Cells.Item(113, 13).Value = -352.875
$ws.Cells.Item(113, 14).Value = -9241.3333
$ws.Cells.Item(129, 8).Value = 16668392
$ws.Cells.Item(129, 9).Value = 31251082
$ws.Cells.Item(129, 10).Value = 2459.4285
$ws.Cells.Item(129, 11).Value = 93753246
$ws.Cells.Item(129, 12).Value = 7378.2855
$ws.Cells.Item(129, 13).Value = -93748246
$ws.Cells.Item(129, 14).Value = -17378.2855
$ws.Cells.Item(132, 8).Value = 4002333
$ws.Cells.Item(132, 9).Value = 4653222
$ws.Cells.Item(132, 10).Value = 4013.8572
$ws.Cells.Item(132, 11).Value = 13959666
$ws.Cells.Item(132, 12).Value = 12041.5716
$ws.Cells.Item(132, 13).Value = -13957136
$ws.Cells.Item(132, 14).Value = -17101.5716
$ws.Cells.Item(136, 8).Value = 30440.77
$ws.Cells.Item(136, 10).Value = 30440.77
$ws.Cells.Item(136, 12).Value = 30440.77
$ws.Cells.Item(136, 14).Value = -40640.77
$ws.Cells.Item(137, 8).Value = 2594.1702
$ws.Cells.Item(137, 9).Value = 2632.5144
$ws.Cells.Item(137, 10).Value = 2482.3333
$ws.Cells.Item(137, 11).Value = 7897.5432
$ws.Cells.Item(137, 12).Value = 7446.999899999999
$ws.Cells.Item(137, 13).Value = -5347.5432
$ws.Cells.Item(137, 14).Value = -12546.9999
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 20836386
$ws.Cells.Item(2, 9).Value = 35716816
$ws.Cells.Item(2, 10).Value = 3782.4
$ws.Cells.Item(2, 11).Value = 35716816
$ws.Cells.Item(2, 12).Value = 3782.4
$ws.Cells.Item(2, 13).Value = -35716703
$ws.Cells.Item(2, 14).Value = -4008.4
$ws.Cells.Item(33, 8).Value = 15600
$ws.Cells.Item(33, 9).Value = 2000
$ws.Cells.Item(33, 10).Value = 19000
$ws.Cells.Item(33, 11).Value = 2000
$ws.Cells.Item(33, 12).Value = 19000
$ws.Cells.Item(33, 13).Value = -1671
$ws.Cells.Item(33, 14).Value = -19658
$ws.Cells.Item(36, 8).Value = 36618.145
$ws.Cells.Item(36, 9).Value = 11560
$ws.Cells.Item(36, 11).Value = 11560
$ws.Cells.Item(36, 13).Value = -11214
$ws.Cells.Item(45, 8).Value = 1357.091
$ws.Cells.Item(45, 9).Value = 1047.8
$ws.Cells.Item(45, 11).Value = 1047.8
$ws.Cells.Item(45, 13).Value = -670.8
$ws.Cells.Item(61, 8).Value = 3004.4348
$ws.Cells.Item(61, 9).Value = 1007.2857
$ws.Cells.Item(61, 11).Value = 1007.2857
$ws.Cells.Item(61, 13).Value = -795.2857
$ws.Cells.Item(74, 8).Value = 524.04346
$ws.Cells.Item(74, 9).Value = 479.2093
$ws.Cells.Item(74, 11).Value = 479.2093
$ws.Cells.Item(74, 13).Value = 394.7907
$ws.Cells.Item(77, 8).Value = 524.04346
$ws.Cells.Item(77, 9).Value = 479.2093
$ws.Cells.Item(77, 11).Value = 2396.0465
$ws.Cells.Item(77, 13).Value = 1971.9535
$ws.Cells.Item(116, 8).Value = 20836386
$ws.Cells.Item(116, 9).Value = 35716816
$ws.Cells.Item(116, 10).Value = 3782.4
$ws.Cells.Item(116, 11).Value = 35716816
$ws.Cells.Item(116, 12).Value = 3782.4
$ws.Cells.Item(116, 13).Value = -35714522
$ws.Cells.Item(116, 14).Value = -8370.4
$ws.Cells.Item(132, 8).Value = 1543.117
$ws.Cells.Item(132, 9).Value = 1031.9
$ws.Cells.Item(132, 10).Value = 3347.4119
$ws.Cells.Item(132, 11).Value = 3095.7
$ws.Cells.Item(132, 12).Value = 10042.2357
$ws.Cells.Item(132, 13).Value = -565.7000000000003
$ws.Cells.Item(132, 14).Value = -15102.2357
$ws.Cells.Item(136, 8).Value = 3004.4348
$ws.Cells.Item(136, 9).Value = 1007.2857
$ws.Cells.Item(136, 11).Value = 3021.8571
$ws.Cells.Item(136, 13).Value = -471.8571000000002
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 20836386
$ws.Cells.Item(3, 9).Value = 35716816
$ws.Cells.Item(3, 10).Value = 3782.4
$ws.Cells.Item(3, 11).Value = 35716816
$ws.Cells.Item(3, 12).Value = 3782.4
$ws.Cells.Item(3, 13).Value = -35716702
$ws.Cells.Item(3, 14).Value = -4010.4
$ws.Cells.Item(20, 8).Value = 6552.4
$ws.Cells.Item(20, 9).Value = 8144.1113
$ws.Cells.Item(20, 10).Value = 4164.8335
$ws.Cells.Item(20, 11).Value = 8144.1113
$ws.Cells.Item(20, 12).Value = 4164.8335
$ws.Cells.Item(20, 13).Value = -7897.1113
$ws.Cells.Item(20, 14).Value = -4658.8335
$ws.Cells.Item(33, 8).Value = 0
$ws.Cells.Item(33, 9).Value = 0
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 11).Value = 0
$ws.Cells.Item(33, 12).Value = 0
$ws.Cells.Item(33, 13).ClearContents()
$ws.Cells.Item(33, 14).ClearContents()
$ws.Cells.Item(133, 8).Value = 40000
$ws.Cells.Item(133, 10).Value = 40000
$ws.Cells.Item(133, 12).Value = 40000
$ws.Cells.Item(133, 14).Value = -50120
$ws.Cells.Item(134, 8).Value = 1733.9811
$ws.Cells.Item(134, 9).Value = 1172.4348
$ws.Cells.Item(134, 10).Value = 5424.143
$ws.Cells.Item(134, 11).Value = 3517.3044
$ws.Cells.Item(134, 12).Value = 16272.429
$ws.Cells.Item(134, 13).Value = -982.3044
$ws.Cells.Item(134, 14).Value = -21342.429
$ws.Cells.Item(135, 8).Value = 29950
$ws.Cells.Item(135, 10).Value = 29950
$ws.Cells.Item(135, 12).Value = 29950
$ws.Cells.Item(135, 14).Value = -40090
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 3360
$ws.Cells.Item(22, 9).Value = 1000
$ws.Cells.Item(22, 10).Value = 3950
$ws.Cells.Item(22, 11).Value = 1000
$ws.Cells.Item(22, 12).Value = 3950
$ws.Cells.Item(22, 13).Value = -650
$ws.Cells.Item(22, 14).Value = -4650
$ws.Cells.Item(31, 8).Value = 2025.2375
$ws.Cells.Item(31, 9).Value = 1300.2245
$ws.Cells.Item(31, 10).Value = 3171.2258
$ws.Cells.Item(31, 11).Value = 1300.2245
$ws.Cells.Item(31, 12).Value = 3171.2258
$ws.Cells.Item(31, 13).Value = -1005.2245
$ws.Cells.Item(31, 14).Value = -3761.2258
$ws.Cells.Item(34, 8).Value = 2025.2375
$ws.Cells.Item(34, 9).Value = 1300.2245
$ws.Cells.Item(34, 10).Value = 3171.2258
$ws.Cells.Item(34, 11).Value = 1300.2245
$ws.Cells.Item(34, 12).Value = 3171.2258
$ws.Cells.Item(34, 13).Value = -1098.2245
$ws.Cells.Item(34, 14).Value = -3575.2258
$ws.Cells.Item(58, 8).Value = 7144735
$ws.Cells.Item(58, 9).Value = 891.2037
$ws.Cells.Item(58, 10).Value = 31255208
$ws.Cells.Item(58, 11).Value = 891.2037
$ws.Cells.Item(58, 12).Value = 31255208
$ws.Cells.Item(58, 13).Value = -688.2037
$ws.Cells.Item(58, 14).Value = -31255614
$ws.Cells.Item(132, 8).Value = 1563.7322
$ws.Cells.Item(132, 9).Value = 1152.7805
$ws.Cells.Item(132, 10).Value = 2687
$ws.Cells.Item(132, 11).Value = 3458.3415
$ws.Cells.Item(132, 12).Value = 8061
$ws.Cells.Item(132, 13).Value = -928.3415000000005
$ws.Cells.Item(132, 14).Value = -13121
$ws.Cells.Item(134, 8).Value = 1778.3928
$ws.Cells.Item(134, 9).Value = 720.913
$ws.Cells.Item(134, 11).Value = 2162.739
$ws.Cells.Item(134, 13).Value = 372.261
$ws.Cells.Item(136, 8).Value = 7144735
$ws.Cells.Item(136, 9).Value = 891.2037
$ws.Cells.Item(136, 10).Value = 31255208
$ws.Cells.Item(136, 11).Value = 2673.6111
$ws.Cells.Item(136, 12).Value = 93765624
$ws.Cells.Item(136, 13).Value = -123.6111000000001
$ws.Cells.Item(136, 14).Value = -93770724
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 1031.7693
$ws.Cells.Item(107, 9).Value = 208.33333
$ws.Cells.Item(107, 10).Value = 1278.8
$ws.Cells.Item(107, 11).Value = 624.99999
$ws.Cells.Item(107, 12).Value = 3836.4
$ws.Cells.Item(107, 13).Value = 1295.00001
$ws.Cells.Item(107, 14).Value = -7676.4
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2184.0962
$ws.Cells.Item(132, 9).Value = 1701.6552
$ws.Cells.Item(132, 10).Value = 2792.3914
$ws.Cells.Item(132, 11).Value = 5104.9656
$ws.Cells.Item(132, 12).Value = 8377.174199999999
$ws.Cells.Item(132, 13).Value = -2574.9656
$ws.Cells.Item(132, 14).Value = -13437.1742
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 90910630
$ws.Cells.Item(22, 9).Value = 250000590
$ws.Cells.Item(22, 10).Value = 2092
$ws.Cells.Item(22, 11).Value = 250000590
$ws.Cells.Item(22, 12).Value = 2092
$ws.Cells.Item(22, 13).Value = -250000295
$ws.Cells.Item(22, 14).Value = -2682
$ws.Cells.Item(27, 8).Value = 90910630
$ws.Cells.Item(27, 9).Value = 250000590
$ws.Cells.Item(27, 10).Value = 2092
$ws.Cells.Item(27, 11).Value = 250000590
$ws.Cells.Item(27, 12).Value = 2092
$ws.Cells.Item(27, 13).Value = -250000483
$ws.Cells.Item(27, 14).Value = -2306
$ws.Cells.Item(46, 8).Value = 2370.1
$ws.Cells.Item(46, 9).Value = 560.2
$ws.Cells.Item(46, 10).Value = 4180
$ws.Cells.Item(46, 11).Value = 560.2
$ws.Cells.Item(46, 12).Value = 4180
$ws.Cells.Item(46, 13).Value = -372.2
$ws.Cells.Item(46, 14).Value = -4556
$ws.Cells.Item(68, 8).Value = 2338.577
$ws.Cells.Item(68, 9).Value = 1083.3334
$ws.Cells.Item(68, 10).Value = 5162.875
$ws.Cells.Item(68, 11).Value = 1083.3334
$ws.Cells.Item(68, 12).Value = 5162.875
$ws.Cells.Item(68, 13).Value = -334.3334
$ws.Cells.Item(68, 14).Value = -6660.875
$ws.Cells.Item(71, 8).Value = 2338.577
$ws.Cells.Item(71, 9).Value = 1083.3334
$ws.Cells.Item(71, 10).Value = 5162.875
$ws.Cells.Item(71, 11).Value = 5416.666999999999
$ws.Cells.Item(71, 12).Value = 25814.375
$ws.Cells.Item(71, 13).Value = -1672.666999999999
$ws.Cells.Item(71, 14).Value = -33302.375
$ws.Cells.Item(132, 8).Value = 1647.3898
$ws.Cells.Item(132, 9).Value = 965.9149
$ws.Cells.Item(132, 11).Value = 2897.7447
$ws.Cells.Item(132, 13).Value = -367.7447000000002
$ws.Cells.Item(136, 8).Value = 1756.7709
$ws.Cells.Item(136, 9).Value = 1333.175
$ws.Cells.Item(136, 10).Value = 4164.8335
$ws.Cells.Item(136, 11).Value = 3999.525
$ws.Cells.Item(136, 13).Value = -1449.525
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 279514.9
$ws.Cells.Item(122, 9).Value = 313926.22
$ws.Cells.Item(122, 11).Value = 941778.6599999999
$ws.Cells.Item(122, 13).Value = -939328.6599999999
$ws.Cells.Item(132, 8).Value = 11729.019
$ws.Cells.Item(132, 9).Value = 2052.2092
$ws.Cells.Item(132, 10).Value = 53339.3
$ws.Cells.Item(132, 11).Value = 6156.6276
$ws.Cells.Item(132, 12).Value = 160017.9
$ws.Cells.Item(132, 13).Value = -3626.6276
$ws.Cells.Item(132, 14).Value = -165077.9
$ws.Cells.Item(136, 8).Value = 977.4314000000001
$ws.Cells.Item(136, 9).Value = 677.1795
$ws.Cells.Item(136, 10).Value = 1953.25
$ws.Cells.Item(136, 11).Value = 2031.5385
$ws.Cells.Item(136, 12).Value = 5859.75
$ws.Cells.Item(136, 14).Value = -10959.75
